# Script 1 - atualização automática de dados (Execução: 22)
# Appends a new data row (year 2023) to the "g1.2" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g1.2")

$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value = 2023
$ws.Cells.Item($newRow, 2).Value = 3.241657824791806
$ws.Cells.Item($newRow, 3).Value = 2.867008788862638
$ws.Cells.Item($newRow, 4).Value = 3.118144130554446
